$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eixos")

# "Domínios" -> "Dominios" (header cell B1 on sheet "Eixos")
$ws.Range("B1").Value = "Dominios"

# Move the active selection on the sheet from B3 to B1
$ws.Range("B1").Select()
